# Update Name of Algo
# Applies the updated result values produced by the KNN imputation algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value  = -7.904000000000001
$ws.Range("D7").Value  = -7.35
$ws.Range("C8").Value  = -12.672
$ws.Range("A12").Value = -21.882
$ws.Range("C12").Value = -13.002
$ws.Range("C14").Value = -11.675
$ws.Range("D19").Value = -7.764000000000001
$ws.Range("E19").Value = 12.965
$ws.Range("D21").Value = -7.597999999999999
$ws.Range("C22").Value = -12.846
$ws.Range("D24").Value = -7.637
